$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 5224.75
$ws.Range("I34").Value = 3359.6
$ws.Range("J34").Value = 8333.333000000001
$ws.Range("K34").Value = 3359.6
$ws.Range("L34").Value = 8333.333000000001
$ws.Range("M34").Value = -3156.6
$ws.Range("N34").Value = -8739.333000000001
$ws.Range("H36").Value = 5224.75
$ws.Range("I36").Value = 3359.6
$ws.Range("J36").Value = 8333.333000000001
$ws.Range("K36").Value = 3359.6
$ws.Range("L36").Value = 8333.333000000001
$ws.Range("M36").Value = -2644.6
$ws.Range("N36").Value = -9763.333000000001
$ws.Range("H98").Value = 7006.7646
$ws.Range("I98").Value = 4587.6665
$ws.Range("J98").Value = 25150
$ws.Range("K98").Value = 4587.6665
$ws.Range("L98").Value = 25150
$ws.Range("M98").Value = -3089.6665
$ws.Range("N98").Value = -28146
$ws.Range("H106").Value = 6944.857
$ws.Range("I106").Value = 7228.2104
$ws.Range("K106").Value = 7228.2104
$ws.Range("M106").Value = -6597.2104
$ws.Range("H122").Value = 7006.7646
$ws.Range("I122").Value = 4587.6665
$ws.Range("J122").Value = 25150
$ws.Range("K122").Value = 13762.9995
$ws.Range("L122").Value = 75450
$ws.Range("M122").Value = -11312.9995
$ws.Range("N122").Value = -80350
$ws.Range("H125").Value = 1684.2222
$ws.Range("I125").Value = 1605.5
$ws.Range("K125").Value = 14449.5
$ws.Range("M125").Value = -11989.5
$ws.Range("H137").Value = 1415.4117
$ws.Range("I137").Value = 975.53845
$ws.Range("J137").Value = 2845
$ws.Range("K137").Value = 2926.61535
$ws.Range("L137").Value = 8535
$ws.Range("M137").Value = -376.61535
$ws.Range("N137").Value = -13635

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13623.875
$ws.Range("I2").Value = 899.5
$ws.Range("K2").Value = 899.5
$ws.Range("M2").Value = -786.5
$ws.Range("H45").Value = 1168.4615
$ws.Range("I45").Value = 1169
$ws.Range("J45").Value = 1166.6666
$ws.Range("K45").Value = 1169
$ws.Range("L45").Value = 1166.6666
$ws.Range("M45").Value = -792
$ws.Range("N45").Value = -1920.6666
$ws.Range("H74").Value = 1153.4333
$ws.Range("I74").Value = 757.9474
$ws.Range("J74").Value = 1836.5454
$ws.Range("K74").Value = 757.9474
$ws.Range("L74").Value = 1836.5454
$ws.Range("M74").Value = 116.0526
$ws.Range("N74").Value = -3584.5454
$ws.Range("H77").Value = 1153.4333
$ws.Range("I77").Value = 757.9474
$ws.Range("J77").Value = 1836.5454
$ws.Range("K77").Value = 3789.737
$ws.Range("L77").Value = 9182.726999999999
$ws.Range("M77").Value = 578.2629999999999
$ws.Range("N77").Value = -17918.727
$ws.Range("H116").Value = 13623.875
$ws.Range("I116").Value = 899.5
$ws.Range("K116").Value = 899.5
$ws.Range("M116").Value = 1394.5
$ws.Range("H122").Value = 1283.1666
$ws.Range("I122").Value = 1079.8
$ws.Range("K122").Value = 3239.4
$ws.Range("M122").Value = -789.3999999999996
$ws.Range("H132").Value = 2006
$ws.Range("I132").Value = 1818.9375
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 5456.8125
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -2926.8125
$ws.Range("N132").Value = -20057
$ws.Range("H133").Value = 31260
$ws.Range("J133").Value = 31260
$ws.Range("L133").Value = 31260
$ws.Range("N133").Value = -36320

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13623.875
$ws.Range("I3").Value = 899.5
$ws.Range("K3").Value = 899.5
$ws.Range("M3").Value = -785.5
$ws.Range("H20").Value = 1448.9333
$ws.Range("I20").Value = 1129.5
$ws.Range("K20").Value = 1129.5
$ws.Range("M20").Value = -882.5
$ws.Range("H80").Value = 670.8946999999999
$ws.Range("I80").Value = 441
$ws.Range("K80").Value = 441
$ws.Range("M80").Value = 557
$ws.Range("H83").Value = 670.8946999999999
$ws.Range("I83").Value = 441
$ws.Range("K83").Value = 2205
$ws.Range("M83").Value = 2787
$ws.Range("H134").Value = 3854.3257
$ws.Range("I134").Value = 1038.875
$ws.Range("J134").Value = 12044.728
$ws.Range("K134").Value = 3116.625
$ws.Range("L134").Value = 36134.18399999999
$ws.Range("M134").Value = -581.625
$ws.Range("N134").Value = -41204.18399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 10278.571
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 10278.571
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 10278.571
$ws.Range("M23").Value = $null
$ws.Range("N23").Value = -10758.571
$ws.Range("H27").Value = 10278.571
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 10278.571
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 10278.571
$ws.Range("M27").Value = $null
$ws.Range("N27").Value = -10662.571
$ws.Range("H31").Value = 1207
$ws.Range("I31").Value = 1189.7646
$ws.Range("J31").Value = 1500
$ws.Range("K31").Value = 1189.7646
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = -894.7646
$ws.Range("N31").Value = -2090
$ws.Range("H34").Value = 1207
$ws.Range("I34").Value = 1189.7646
$ws.Range("J34").Value = 1500
$ws.Range("K34").Value = 1189.7646
$ws.Range("L34").Value = 1500
$ws.Range("M34").Value = -987.7646
$ws.Range("N34").Value = -1904
$ws.Range("H58").Value = 1785
$ws.Range("I58").Value = 1596.6666
$ws.Range("J58").Value = 2350
$ws.Range("K58").Value = 1596.6666
$ws.Range("L58").Value = 2350
$ws.Range("M58").Value = -1393.6666
$ws.Range("N58").Value = -2756
$ws.Range("H122").Value = 764.7619
$ws.Range("I122").Value = 691.3125
$ws.Range("K122").Value = 2073.9375
$ws.Range("M122").Value = 376.0625
$ws.Range("H136").Value = 1785
$ws.Range("I136").Value = 1596.6666
$ws.Range("J136").Value = 2350
$ws.Range("K136").Value = 4789.9998
$ws.Range("L136").Value = 7050
$ws.Range("M136").Value = -2239.9998
$ws.Range("N136").Value = -12150

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null
$ws.Range("H131").Value = 20003026
$ws.Range("J131").Value = 3447.4187
$ws.Range("L131").Value = 10342.2561
$ws.Range("N131").Value = -20422.2561
$ws.Range("H137").Value = 2061.45
$ws.Range("I137").Value = 952.5
$ws.Range("K137").Value = 2857.5
$ws.Range("M137").Value = 2242.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 762
$ws.Range("I97").Value = 762
$ws.Range("K97").Value = 762
$ws.Range("M97").Value = -266
$ws.Range("H113").Value = 1228.1177
$ws.Range("I113").Value = 1199.875
$ws.Range("K113").Value = 1199.875
$ws.Range("M113").Value = 970.125
$ws.Range("H122").Value = 2444.111
$ws.Range("I122").Value = 2499.625
$ws.Range("K122").Value = 7498.875
$ws.Range("M122").Value = -5048.875
$ws.Range("H132").Value = 1723.0322
$ws.Range("I132").Value = 1369.3334
$ws.Range("J132").Value = 4110.5
$ws.Range("K132").Value = 4108.0002
$ws.Range("L132").Value = 12331.5
$ws.Range("M132").Value = -1578.0002
$ws.Range("N132").Value = -17391.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1428.3077
$ws.Range("I7").Value = 1192.55
$ws.Range("J7").Value = 2214.1667
$ws.Range("K7").Value = 1192.55
$ws.Range("L7").Value = 2214.1667
$ws.Range("M7").Value = -1080.55
$ws.Range("N7").Value = -2438.1667
$ws.Range("H46").Value = 2085.4285
$ws.Range("I46").Value = 933
$ws.Range("J46").Value = 2949.75
$ws.Range("K46").Value = 933
$ws.Range("L46").Value = 2949.75
$ws.Range("M46").Value = -745
$ws.Range("N46").Value = -3325.75
$ws.Range("H55").Value = 205.03448
$ws.Range("I55").Value = 159.57143
$ws.Range("K55").Value = 159.57143
$ws.Range("M55").Value = 13.42857000000001
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = $null
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = $null
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null
$ws.Range("H126").Value = 1428.3077
$ws.Range("I126").Value = 1192.55
$ws.Range("J126").Value = 2214.1667
$ws.Range("K126").Value = 3577.65
$ws.Range("L126").Value = 6642.500100000001
$ws.Range("M126").Value = -1107.65
$ws.Range("N126").Value = -11582.5001
$ws.Range("H136").Value = 1407.7646
$ws.Range("I136").Value = 1339.1875
$ws.Range("J136").Value = 2505
$ws.Range("K136").Value = 4017.5625
$ws.Range("L136").Value = 7515
$ws.Range("M136").Value = -1467.5625
$ws.Range("N136").Value = -12615
$ws.Range("H140").Value = 83250
$ws.Range("J140").Value = 83250
$ws.Range("L140").Value = 83250
$ws.Range("N140").Value = -93610

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 57274150
$ws.Range("I122").Value = 78751304
$ws.Range("J122").Value = 1731.6666
$ws.Range("K122").Value = 236253912
$ws.Range("L122").Value = 5194.9998
$ws.Range("M122").Value = -236251462
$ws.Range("N122").Value = -10094.9998
$ws.Range("H126").Value = 46297136
$ws.Range("I126").Value = 52910900
$ws.Range("J126").Value = 793.3333
$ws.Range("K126").Value = 158732700
$ws.Range("L126").Value = 2379.9999
$ws.Range("M126").Value = -158730230
$ws.Range("N126").Value = -7319.9999
